$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.858.82"
$ws.Range("E2").Value = "  +4.79%  "

$ws.Range("D3").Value = "2.265.17"
$ws.Range("E3").Value = "  +2.06%  "

$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "301.85"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "92.00"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.93%  "

$ws.Range("E7").Value = "  +3.35%  "

$ws.Range("E8").Value = "  -0.08%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.484"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.92%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "54.28"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.77%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "32.22"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.79%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0797"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.31%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("E14").Value = "  +3.53%  "

$ws.Range("D15").Value = "2.618.18"
$ws.Range("E15").Value = "  +2.23%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.15"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.65%  "

$ws.Range("D17").Value = "2.272.26"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("E18").Value = "  +3.42%  "

$ws.Range("D19").Value = "41.753.26"
$ws.Range("E19").Value = "  +4.73%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.07"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +8.62%  "

$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +2.01%  "

$ws.Range("E22").Value = "  +3.54%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "66.95"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.15%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "241.62"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.78%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  +3.85%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "23.86"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.53%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.62"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.25%  "

$ws.Range("E30").Value = "  -12.06%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "159.44"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "33.74"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +6.07%  "

$ws.Range("E33").Value = "  -0.02%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.16"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.91%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0745"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.24%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.08"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("E37").Value = "  +2.03%  "

$ws.Range("E38").Value = "  +5.25%  "

$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("E40").Value = "  +8.86%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.81"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.91"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.61%  "

$ws.Range("D43").Value = "2.073.56"
$ws.Range("E43").Value = "  -0.56%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "19.52"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +9.25%  "

$ws.Range("E45").Value = "  +3.18%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.16"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.86%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.51"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.56%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.14"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.19%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "51.70"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.60%  "
